$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1962.8889
$ws.Range("I19").Value = 1442.8572
$ws.Range("J19").Value = 2293.818
$ws.Range("K19").Value = 1442.8572
$ws.Range("L19").Value = 2293.818
$ws.Range("M19").Value = -1267.8572
$ws.Range("N19").Value = -2643.818
# Row 51
$ws.Range("H51").Value = 4204.1665
$ws.Range("I51").Value = 3125
$ws.Range("J51").Value = 4420
$ws.Range("K51").Value = 3125
$ws.Range("L51").Value = 4420
$ws.Range("M51").Value = -2641
$ws.Range("N51").Value = -5388
# Row 123
$ws.Range("H123").Value = 40057.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 40057.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 40057.5
$ws.Range("N123").Value = -49857.5
# Row 135
$ws.Range("H135").Value = 12928348
$ws.Range("I135").Value = 427.79166
$ws.Range("J135").Value = 31179530
$ws.Range("K135").Value = 3850.12494
$ws.Range("L135").Value = 280615770
$ws.Range("M135").Value = -1315.12494
# Row 137
$ws.Range("H137").Value = 26317174
$ws.Range("I137").Value = 1285.76
$ws.Range("J137").Value = 76924650
$ws.Range("K137").Value = 3857.28
$ws.Range("L137").Value = 230773950
$ws.Range("M137").Value = -1307.28
$ws.Range("N137").Value = -230779050

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 774936.9399999999
$ws.Range("I2").Value = 665
$ws.Range("J2").Value = 2452526
$ws.Range("K2").Value = 665
$ws.Range("L2").Value = 2452526
$ws.Range("M2").Value = -552
$ws.Range("N2").Value = -2452752
# Row 63
$ws.Range("H63").Value = 2815.182
$ws.Range("I63").Value = 2040.7778
$ws.Range("J63").Value = 6300
$ws.Range("K63").Value = 2040.7778
$ws.Range("L63").Value = 6300
$ws.Range("M63").Value = -1354.7778
$ws.Range("N63").Value = -7672
# Row 66
$ws.Range("H66").Value = 2815.182
$ws.Range("I66").Value = 2040.7778
$ws.Range("J66").Value = 6300
$ws.Range("K66").Value = 10203.889
$ws.Range("L66").Value = 31500
$ws.Range("M66").Value = -6771.889000000001
$ws.Range("N66").Value = -38364
# Row 116
$ws.Range("H116").Value = 774936.9399999999
$ws.Range("I116").Value = 665
$ws.Range("J116").Value = 2452526
$ws.Range("K116").Value = 665
$ws.Range("L116").Value = 2452526
$ws.Range("M116").Value = 1629
$ws.Range("N116").Value = -2457114
# Row 122
$ws.Range("H122").Value = 1336.7576
$ws.Range("I122").Value = 1301.7727
$ws.Range("J122").Value = 1406.7273
$ws.Range("K122").Value = 3905.3181
$ws.Range("L122").Value = 4220.1819
$ws.Range("M122").Value = -1455.3181
$ws.Range("N122").Value = -9120.1819
# Row 132
$ws.Range("H132").Value = 671048.6
$ws.Range("I132").Value = 822648.5600000001
$ws.Range("J132").Value = 92212.45
$ws.Range("K132").Value = 2467945.68
$ws.Range("L132").Value = 276637.35
$ws.Range("M132").Value = -2465415.68
$ws.Range("N132").Value = -281697.35

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 774936.9399999999
$ws.Range("I3").Value = 665
$ws.Range("J3").Value = 2452526
$ws.Range("K3").Value = 665
$ws.Range("L3").Value = 2452526
$ws.Range("M3").Value = -551
$ws.Range("N3").Value = -2452754
# Row 107
$ws.Range("H107").Value = 705370.9
$ws.Range("I107").Value = 1084417.1
$ws.Range("J107").Value = 1428
$ws.Range("K107").Value = 1084417.1
$ws.Range("L107").Value = 1428
$ws.Range("M107").Value = -1082497.1
$ws.Range("N107").Value = -5268

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 332969.66
$ws.Range("I31").Value = 1086.4375
$ws.Range("J31").Value = 815708.9
$ws.Range("K31").Value = 1086.4375
$ws.Range("L31").Value = 815708.9
$ws.Range("M31").Value = -791.4375
$ws.Range("N31").Value = -816298.9
# Row 34
$ws.Range("H34").Value = 332969.66
$ws.Range("I34").Value = 1086.4375
$ws.Range("J34").Value = 815708.9
$ws.Range("K34").Value = 1086.4375
$ws.Range("L34").Value = 815708.9
$ws.Range("M34").Value = -884.4375
$ws.Range("N34").Value = -816112.9
# Row 132
$ws.Range("H132").Value = 1849.591
$ws.Range("I132").Value = 1746.7
$ws.Range("J132").Value = 2878.5
$ws.Range("K132").Value = 5240.1
$ws.Range("L132").Value = 8635.5
$ws.Range("M132").Value = -2710.1
$ws.Range("N132").Value = -13695.5
# Row 134
$ws.Range("H134").Value = 1019.3051
$ws.Range("I134").Value = 1074.4423
$ws.Range("J134").Value = 609.7143
$ws.Range("K134").Value = 3223.3269
$ws.Range("L134").Value = 1829.1429
$ws.Range("M134").Value = -688.3269
$ws.Range("N134").Value = -6899.1429

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
# Row 76
$ws.Range("H76").Value = 3163
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3163
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 9489
$ws.Range("N76").Value = -10255
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 3163
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3163
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 9489
$ws.Range("N79").Value = -12141
$ws.Range("M79").ClearContents()
# Row 140
$ws.Range("H140").Value = 1875
$ws.Range("I140").Value = 1065
$ws.Range("J140").Value = 3900
$ws.Range("K140").Value = 3195
$ws.Range("L140").Value = 11700
$ws.Range("M140").Value = 1985
$ws.Range("N140").Value = -22060
# Row 141
$ws.Range("H141").Value = 2000000
$ws.Range("I141").Value = 2000000
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6000000
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5994820
$ws.Range("N141").ClearContents()

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1635.0178
$ws.Range("I132").Value = 1608.1904
$ws.Range("J132").Value = 1715.5
$ws.Range("K132").Value = 4824.5712
$ws.Range("L132").Value = 5146.5
$ws.Range("M132").Value = -2294.5712
$ws.Range("N132").Value = -10206.5

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1825.2084
$ws.Range("I68").Value = 1700.125
$ws.Range("J68").Value = 2075.375
$ws.Range("K68").Value = 1700.125
$ws.Range("L68").Value = 2075.375
$ws.Range("M68").Value = -951.125
$ws.Range("N68").Value = -3573.375
# Row 71
$ws.Range("H71").Value = 1825.2084
$ws.Range("I71").Value = 1700.125
$ws.Range("J71").Value = 2075.375
$ws.Range("K71").Value = 8500.625
$ws.Range("L71").Value = 10376.875
$ws.Range("M71").Value = -4756.625
$ws.Range("N71").Value = -17864.875
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
# Row 132
$ws.Range("H132").Value = 4508.3237
$ws.Range("I132").Value = 4875.0645
$ws.Range("J132").Value = 718.6667
$ws.Range("K132").Value = 14625.1935
$ws.Range("L132").Value = 2156.0001
$ws.Range("M132").Value = -12095.1935
$ws.Range("N132").Value = -7216.0001
# Row 136
$ws.Range("H136").Value = 1470.2858
$ws.Range("I136").Value = 1140.7407
$ws.Range("J136").Value = 2582.5
$ws.Range("K136").Value = 3422.2221
$ws.Range("L136").Value = 7747.5
$ws.Range("M136").Value = -872.2221
$ws.Range("N136").Value = -12847.5

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3385.0925
$ws.Range("I132").Value = 3761.9788
$ws.Range("J132").Value = 854.5714
$ws.Range("K132").Value = 11285.9364
$ws.Range("L132").Value = 2563.7142
$ws.Range("M132").Value = -8755.936399999999
$ws.Range("N132").Value = -7623.7142
# Row 136
$ws.Range("H136").Value = 4216.755
$ws.Range("I136").Value = 4682.4683
$ws.Range("J136").Value = 568.6667
$ws.Range("K136").Value = 14047.4049
$ws.Range("L136").Value = 1706.0001
$ws.Range("M136").Value = -11497.4049
$ws.Range("N136").Value = -6806.0001
